$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 54, shifting existing rows 54-145 down to 55-146
$ws.Rows(54).Insert()

# Populate the newly inserted row 54 with the new weekly record
$ws.Range("A54").Value = 5
$ws.Range("B54").Value = "Macroferia Regional de Talca"
$ws.Range("C54").Value = "Maule"
$ws.Range("D54").Value = 45219
$ws.Range("E54").Value = 7
$ws.Range("F54").Value = 100112022
$ws.Range("G54").Value = "Arveja Verde"
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 200
$ws.Range("K54").Value = 27000
$ws.Range("L54").Value = 27000
$ws.Range("M54").Value = 27000
$ws.Range("N54").Value = "$/saco 25 kilos"
$ws.Range("O54").Value = "Región del Maule"
$ws.Range("P54").Value = 1080
$ws.Range("Q54").Value = 25
$ws.Range("R54").Value = "Hortaliza"
